$wb = $excel.ActiveWorkbook

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1114.4584
$ws.Range("I129").Value = 444
$ws.Range("J129").Value = 1143.6086
$ws.Range("K129").Value = 1332
$ws.Range("L129").Value = 3430.8258
$ws.Range("M129").Value = 3668
$ws.Range("N129").Value = -13430.8258

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 557226.4399999999
$ws.Range("I132").Value = 675054.2
$ws.Range("J132").Value = 27001.5
$ws.Range("K132").Value = 2025162.6
$ws.Range("L132").Value = 81004.5
$ws.Range("M132").Value = -2022632.6
$ws.Range("N132").Value = -86064.5

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1060
$ws.Range("I45").Value = 950
$ws.Range("K45").Value = 950
$ws.Range("M45").Value = -573

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 14220
$ws.Range("I63").Value = 14220
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 14220
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -13534
$ws.Range("N63").ClearContents()

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 14220
$ws.Range("I66").Value = 14220
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 71100
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -67668
$ws.Range("N66").ClearContents()

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3915.4285
$ws.Range("I88").Value = 3834.3333
$ws.Range("J88").Value = 3976.25
$ws.Range("K88").Value = 3834.3333
$ws.Range("L88").Value = 3976.25
$ws.Range("M88").Value = -3428.3333
$ws.Range("N88").Value = -4788.25

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3915.4285
$ws.Range("I91").Value = 3834.3333
$ws.Range("J91").Value = 3976.25
$ws.Range("K91").Value = 3834.3333
$ws.Range("L91").Value = 3976.25
$ws.Range("M91").Value = -2430.3333
$ws.Range("N91").Value = -6784.25

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 37048544
$ws.Range("I97").Value = 47633696
$ws.Range("J97").Value = 499.5
$ws.Range("K97").Value = 47633696
$ws.Range("L97").Value = 499.5
$ws.Range("M97").Value = -47633200
$ws.Range("N97").Value = -1491.5

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 275
$ws.Range("I22").Value = 275
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 275
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -102
$ws.Range("N22").ClearContents()

# BSM row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 27500
$ws.Range("J100").Value = 27500
$ws.Range("L100").Value = 27500
$ws.Range("N100").Value = -29664

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 336710.53
$ws.Range("I105").Value = 3365.7896
$ws.Range("J105").Value = 912487.8
$ws.Range("K105").Value = 3365.7896
$ws.Range("L105").Value = 912487.8
$ws.Range("M105").Value = -1618.7896
$ws.Range("N105").Value = -915981.8

# BSM row 133
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1987.2667
$ws.Range("I94").Value = 1180
$ws.Range("J94").Value = 2390.9
$ws.Range("K94").Value = 1180
$ws.Range("L94").Value = 2390.9
$ws.Range("M94").Value = -729
$ws.Range("N94").Value = -3292.9

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2585.1025
$ws.Range("I134").Value = 1200.3914
$ws.Range("J134").Value = 4575.625
$ws.Range("K134").Value = 3601.1742
$ws.Range("L134").Value = 13726.875
$ws.Range("M134").Value = -1066.1742
$ws.Range("N134").Value = -18796.875

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5990
$ws.Range("I56").Value = 5990
$ws.Range("K56").Value = 5990
$ws.Range("M56").Value = -5460

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4763147
$ws.Range("J131").Value = 5129526
$ws.Range("L131").Value = 15388578
$ws.Range("N131").Value = -15398658

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6729.8887
$ws.Range("I70").Value = 7973.4546
$ws.Range("K70").Value = 7973.4546
$ws.Range("M70").Value = -7703.4546

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6729.8887
$ws.Range("I73").Value = 7973.4546
$ws.Range("K73").Value = 7973.4546
$ws.Range("M73").Value = -7037.4546

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2961.5
$ws.Range("I80").Value = 2912.6667
$ws.Range("J80").Value = 3042.889
$ws.Range("K80").Value = 2912.6667
$ws.Range("L80").Value = 3042.889
$ws.Range("M80").Value = -1914.6667
$ws.Range("N80").Value = -5038.889

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2961.5
$ws.Range("I83").Value = 2912.6667
$ws.Range("J83").Value = 3042.889
$ws.Range("K83").Value = 14563.3335
$ws.Range("L83").Value = 15214.445
$ws.Range("M83").Value = -9571.333500000001
$ws.Range("N83").Value = -25198.445

# GSM row 116
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 32000
$ws.Range("J116").Value = 32000
$ws.Range("L116").Value = 32000
$ws.Range("N116").Value = -41178

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2400.3333
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2600.5
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2600.5
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -4098.5

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2400.3333
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2600.5
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 13002.5
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -20490.5

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2850.7693
$ws.Range("I82").Value = 3063
$ws.Range("J82").Value = 2756.4443
$ws.Range("K82").Value = 3063
$ws.Range("L82").Value = 2756.4443
$ws.Range("M82").Value = -2702
$ws.Range("N82").Value = -3478.4443

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2850.7693
$ws.Range("I85").Value = 3063
$ws.Range("J85").Value = 2756.4443
$ws.Range("K85").Value = 3063
$ws.Range("L85").Value = 2756.4443
$ws.Range("M85").Value = -1815
$ws.Range("N85").Value = -5252.4443

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1916.6666
$ws.Range("I93").Value = 875
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 875
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = 373
$ws.Range("N93").Value = -6496

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3207.2
$ws.Range("I132").Value = 1854
$ws.Range("J132").Value = 4976.769
$ws.Range("K132").Value = 5562
$ws.Range("L132").Value = 14930.307
$ws.Range("M132").Value = -3032
$ws.Range("N132").Value = -19990.307

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 67845.39999999999
$ws.Range("I122").Value = 84323.414
$ws.Range("K122").Value = 252970.242
$ws.Range("M122").Value = -250520.242

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2040.8387
$ws.Range("I136").Value = 1359.9412
$ws.Range("J136").Value = 2867.6428
$ws.Range("K136").Value = 4079.8236
$ws.Range("L136").Value = 8602.928400000001
$ws.Range("M136").Value = -1529.8236
$ws.Range("N136").Value = -13702.9284
